# "Added paypal 2 request"
#
# testdata.xlsx has a "deleteCustomer" section on the "testdata" sheet
# (A13:A15) whose value cell (A15, under the "id" header in A14) held a
# stale Stripe/PayPal customer id left over from a previous test run.
# This commit swaps that id for the one created by the new "paypal 2"
# request, and leaves the workbook focused on the "testdata" sheet
# (which is where that id was just edited) instead of the
# "cxCreationInvalidKey" sheet that was active before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# Replace the old customer id used by the deleteCustomer test case.
$ws.Range("A15").Value = "cus_OehjJmR5GBWHrF"

# Leave the "testdata" sheet active/selected (it was last touched here),
# with B14 as the active cell, matching where editing left off.
$ws.Activate()
$ws.Range("B14").Select() | Out-Null
